# Extract XLSX common code from conference scheduling example
# This reproduces, via Excel COM interop, the effects of the upstream commit:
#  - "Score view" sheet becomes "Infeasible view" (and its usable-session count is fixed: 60 -> 36)
#  - "Audience type view" / "Audience level view" are pluralised
#  - "Rooms view" row heights shrink (60 -> 45)
#  - "Speakers view" gains many more (wider) columns
#  - All pinned-talk comments get reformatted to the new "total/breakdown" wording
#  - Two new sheets are appended: "Languages view" and a brand-new "Score view"
#    that lists the constraint-match breakdown.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename existing sheets
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Score view").Name = "Infeasible view"
$wb.Worksheets.Item("Audience type view").Name = "Audience types view"
$wb.Worksheets.Item("Audience level view").Name = "Audience levels view"

# ---------------------------------------------------------------------------
# 2. "Infeasible view" (formerly "Score view"): usable sessions 60 -> 36
# ---------------------------------------------------------------------------
$infeasible = $wb.Worksheets.Item("Infeasible view")
$infeasible.Range("E7").Value = 36

# ---------------------------------------------------------------------------
# 3. "Rooms view": row heights 60 -> 45 for rows 3-7
# ---------------------------------------------------------------------------
$roomsView = $wb.Worksheets.Item("Rooms view")
for ($r = 3; $r -le 7; $r++) {
    $roomsView.Rows.Item($r).RowHeight = 45
}

# ---------------------------------------------------------------------------
# 4. "Speakers view": widen/extend the columns (13 bestFit cols -> 25 fixed
#    20-char-wide columns, columns 2..26)
# ---------------------------------------------------------------------------
$speakersView = $wb.Worksheets.Item("Speakers view")
for ($c = 2; $c -le 26; $c++) {
    $speakersView.Columns.Item($c).ColumnWidth = 19.2
}

# ---------------------------------------------------------------------------
# 5. Update the pinned-talk comment text on every sheet that carries one.
# ---------------------------------------------------------------------------
$newCommentText = "S14: Troubleshooting reliable RestEasy`n    Amy Green`nPINNED BY USER`n-1hard total`n    -1hard for 1 Speaker unavailable timeslots`n"

$roomsView.Range("C3").Comment.Text($newCommentText)
$speakersView.Range("C23").Comment.Text($newCommentText)
$wb.Worksheets.Item("Theme tracks view").Range("C3").Comment.Text($newCommentText)
$wb.Worksheets.Item("Audience types view").Range("C3").Comment.Text($newCommentText)
$wb.Worksheets.Item("Audience levels view").Range("C3").Comment.Text($newCommentText)
$wb.Worksheets.Item("Contents view").Range("C3").Comment.Text($newCommentText)

Write-Host "Phase 1 (renames, tweaks, comments) done"
